# Refresh the crypto price/volume snapshot (and correct the Polkadot /
# Chainlink row ordering) per the automated GitHub Actions data pull on
# 2023-05-06. All Price/Volume cells are plain text in this sheet, so
# values that Excel would otherwise reinterpret as numbers (e.g. "1.000"
# or "17.00", which would lose their trailing zeros) are written with a
# leading apostrophe to force literal text storage; values that are
# already unambiguous text (containing extra dots, "%", spaces, etc.) are
# written as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Col = 4; Value = "29.256.28" },
    @{ Row = 2; Col = 5; Value = "  +0.40%  " },
    @{ Row = 3; Col = 4; Value = "1.929.78" },
    @{ Row = 3; Col = 5; Value = "  +1.26%  " },
    @{ Row = 4; Col = 4; Value = "'1.000" },
    @{ Row = 4; Col = 5; Value = "  -0.14%  " },
    @{ Row = 5; Col = 4; Value = "'325.31" },
    @{ Row = 5; Col = 5; Value = "  -0.16%  " },
    @{ Row = 6; Col = 4; Value = "'0.9992" },
    @{ Row = 7; Col = 4; Value = "'0.4615" },
    @{ Row = 7; Col = 5; Value = "  +0.11%  " },
    @{ Row = 8; Col = 4; Value = "'0.3864" },
    @{ Row = 8; Col = 5; Value = "  -0.71%  " },
    @{ Row = 9; Col = 4; Value = "'45.78" },
    @{ Row = 9; Col = 5; Value = "  -1.09%  " },
    @{ Row = 10; Col = 4; Value = "'0.07785" },
    @{ Row = 10; Col = 5; Value = "  -1.21%  " },
    @{ Row = 11; Col = 4; Value = "'0.9703" },
    @{ Row = 11; Col = 5; Value = "  -2.02%  " },
    @{ Row = 12; Col = 4; Value = "'22.55" },
    @{ Row = 12; Col = 5; Value = "  +2.55%  " },
    @{ Row = 13; Col = 4; Value = "1.929.75" },
    @{ Row = 13; Col = 5; Value = "  +2.48%  " },
    @{ Row = 14; Col = 2; Value = "Polkadot" },
    @{ Row = 14; Col = 3; Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot" },
    @{ Row = 14; Col = 4; Value = "'5.770" },
    @{ Row = 14; Col = 5; Value = "  +0.04%  " },
    @{ Row = 15; Col = 2; Value = "Chainlink" },
    @{ Row = 15; Col = 3; Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link" },
    @{ Row = 15; Col = 4; Value = "'7.062" },
    @{ Row = 15; Col = 5; Value = "  +0.31%  " },
    @{ Row = 16; Col = 5; Value = "  +0.29%  " },
    @{ Row = 17; Col = 4; Value = "'86.61" },
    @{ Row = 17; Col = 5; Value = "  -1.72%  " },
    @{ Row = 18; Col = 4; Value = "'1.002" },
    @{ Row = 18; Col = 5; Value = "  -0.12%  " },
    @{ Row = 19; Col = 4; Value = "'0.000009642" },
    @{ Row = 19; Col = 5; Value = "  -3.12%  " },
    @{ Row = 20; Col = 4; Value = "'17.00" },
    @{ Row = 20; Col = 5; Value = "  -0.43%  " },
    @{ Row = 21; Col = 4; Value = "'0.9999" },
    @{ Row = 21; Col = 5; Value = "  -0.19%  " },
    @{ Row = 22; Col = 4; Value = "29.250.26" },
    @{ Row = 23; Col = 4; Value = "'5.462" },
    @{ Row = 23; Col = 5; Value = "  +2.62%  " },
    @{ Row = 24; Col = 4; Value = "'11.03" },
    @{ Row = 24; Col = 5; Value = "  -0.87%  " },
    @{ Row = 25; Col = 4; Value = "2.161.50" },
    @{ Row = 25; Col = 5; Value = "  +1.76%  " },
    @{ Row = 27; Col = 4; Value = "'156.73" },
    @{ Row = 27; Col = 5; Value = "  +0.24%  " },
    @{ Row = 28; Col = 4; Value = "'19.31" },
    @{ Row = 28; Col = 5; Value = "  -0.81%  " },
    @{ Row = 29; Col = 4; Value = "'5.740" },
    @{ Row = 29; Col = 5; Value = "  -2.84%  " },
    @{ Row = 30; Col = 4; Value = "'118.16" },
    @{ Row = 30; Col = 5; Value = "  -0.57%  " },
    @{ Row = 31; Col = 4; Value = "'1.846" },
    @{ Row = 31; Col = 5; Value = "  -1.76%  " },
    @{ Row = 32; Col = 4; Value = "'0.09336" },
    @{ Row = 32; Col = 5; Value = "  -0.22%  " },
    @{ Row = 33; Col = 4; Value = "'0.8586" },
    @{ Row = 33; Col = 5; Value = "  -4.20%  " },
    @{ Row = 34; Col = 4; Value = "'5.156" },
    @{ Row = 34; Col = 5; Value = "  -1.42%  " },
    @{ Row = 35; Col = 4; Value = "'1.302" },
    @{ Row = 35; Col = 5; Value = "  -1.56%  " },
    @{ Row = 36; Col = 4; Value = "'3.073" },
    @{ Row = 36; Col = 5; Value = "  -2.62%  " },
    @{ Row = 37; Col = 4; Value = "'0.05759" },
    @{ Row = 37; Col = 5; Value = "  -0.62%  " },
    @{ Row = 38; Col = 5; Value = "  -1.46%  " },
    @{ Row = 39; Col = 4; Value = "'0.02072" },
    @{ Row = 39; Col = 5; Value = "  -0.69%  " },
    @{ Row = 40; Col = 4; Value = "'7.622" },
    @{ Row = 40; Col = 5; Value = "  -0.72%  " },
    @{ Row = 41; Col = 4; Value = "'0.5637" },
    @{ Row = 42; Col = 4; Value = "'0.000003140" },
    @{ Row = 42; Col = 5; Value = "  +56.69%  " },
    @{ Row = 43; Col = 4; Value = "'0.1773" },
    @{ Row = 43; Col = 5; Value = "  -1.96%  " },
    @{ Row = 44; Col = 4; Value = "'9.352" },
    @{ Row = 44; Col = 5; Value = "  -3.76%  " },
    @{ Row = 45; Col = 4; Value = "'2.709" },
    @{ Row = 45; Col = 5; Value = "  +6.23%  " },
    @{ Row = 46; Col = 4; Value = "'0.5259" },
    @{ Row = 46; Col = 5; Value = "  -1.86%  " },
    @{ Row = 47; Col = 4; Value = "'11.50" },
    @{ Row = 47; Col = 5; Value = "  -3.09%  " },
    @{ Row = 48; Col = 5; Value = "  -2.13%  " },
    @{ Row = 49; Col = 4; Value = "'2.077" },
    @{ Row = 49; Col = 5; Value = "  -4.67%  " },
    @{ Row = 50; Col = 5; Value = "  -1.80%  " },
    @{ Row = 51; Col = 5; Value = "  -1.83%  " }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

Write-Host "Applied $($updates.Count) cell updates"